$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1127.4286
$ws.Range("I19").Value = 878.8
$ws.Range("J19").Value = 1265.5555
$ws.Range("K19").Value = 878.8
$ws.Range("L19").Value = 1265.5555
$ws.Range("M19").Value = -703.8
$ws.Range("N19").Value = -1615.5555
$ws.Range("H40").Value = 1960.48
$ws.Range("I40").Value = 2253.0588
$ws.Range("J40").Value = 1338.75
$ws.Range("K40").Value = 2253.0588
$ws.Range("L40").Value = 1338.75
$ws.Range("M40").Value = -2078.0588
$ws.Range("N40").Value = -1688.75
$ws.Range("H70").Value = 1452.6111
$ws.Range("I70").Value = 1595
$ws.Range("J70").Value = 1434.8125
$ws.Range("K70").Value = 4785
$ws.Range("L70").Value = 4304.4375
$ws.Range("M70").Value = -4515
$ws.Range("N70").Value = -4844.4375
$ws.Range("H73").Value = 1452.6111
$ws.Range("I73").Value = 1595
$ws.Range("J73").Value = 1434.8125
$ws.Range("K73").Value = 4785
$ws.Range("L73").Value = 4304.4375
$ws.Range("M73").Value = -3849
$ws.Range("N73").Value = -6176.4375
$ws.Range("H137").Value = 1729.0847
$ws.Range("I137").Value = 1337.2609
$ws.Range("J137").Value = 3115.5386
$ws.Range("K137").Value = 4011.7827
$ws.Range("L137").Value = 9346.6158
$ws.Range("M137").Value = -1461.7827
$ws.Range("N137").Value = -14446.6158
$ws.Range("H138").Value = 3672.283
$ws.Range("I138").Value = 1442.0555
$ws.Range("J138").Value = 4819.2573
$ws.Range("K138").Value = 4326.166499999999
$ws.Range("L138").Value = 14457.7719
$ws.Range("M138").Value = 813.8335000000006
$ws.Range("N138").Value = -24737.7719

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15428.447
$ws.Range("I32").Value = 16433.445
$ws.Range("J32").Value = 8667.546
$ws.Range("K32").Value = 16433.445
$ws.Range("L32").Value = 8667.546
$ws.Range("M32").Value = -16146.445
$ws.Range("N32").Value = -9241.546
$ws.Range("H63").Value = 4665
$ws.Range("I63").Value = 4502.5
$ws.Range("K63").Value = 4502.5
$ws.Range("M63").Value = -3816.5
$ws.Range("H64").Value = 38091
$ws.Range("J64").Value = 38091
$ws.Range("L64").Value = 38091
$ws.Range("N64").Value = -38587
$ws.Range("H66").Value = 4665
$ws.Range("I66").Value = 4502.5
$ws.Range("K66").Value = 22512.5
$ws.Range("M66").Value = -19080.5
$ws.Range("H67").Value = 38091
$ws.Range("J67").Value = 38091
$ws.Range("L67").Value = 38091
$ws.Range("N67").Value = -39807
$ws.Range("H74").Value = 4865.55
$ws.Range("I74").Value = 3416.745
$ws.Range("J74").Value = 13075.444
$ws.Range("K74").Value = 3416.745
$ws.Range("L74").Value = 13075.444
$ws.Range("M74").Value = -2542.745
$ws.Range("N74").Value = -14823.444
$ws.Range("H77").Value = 4865.55
$ws.Range("I77").Value = 3416.745
$ws.Range("J77").Value = 13075.444
$ws.Range("K77").Value = 17083.725
$ws.Range("L77").Value = 65377.22
$ws.Range("M77").Value = -12715.725
$ws.Range("N77").Value = -74113.22
$ws.Range("H110").Value = 1818.6666
$ws.Range("I110").Value = 1818.6666
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1818.6666
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 226.3334
$ws.Range("N110").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 14949.75
$ws.Range("I7").Value = 9999
$ws.Range("J7").Value = 16600
$ws.Range("K7").Value = 9999
$ws.Range("L7").Value = 16600
$ws.Range("M7").Value = -9886
$ws.Range("N7").Value = -16826
$ws.Range("H12").Value = 5070
$ws.Range("I12").Value = 1240
$ws.Range("J12").Value = 8900
$ws.Range("K12").Value = 1240
$ws.Range("L12").Value = 8900
$ws.Range("M12").Value = -1072
$ws.Range("N12").Value = -9236
$ws.Range("H18").Value = 15133.333
$ws.Range("J18").Value = 15133.333
$ws.Range("L18").Value = 15133.333
$ws.Range("N18").Value = -16191.333
$ws.Range("H22").Value = 484
$ws.Range("I22").Value = 473.33334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 473.33334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -300.33334
$ws.Range("N22").Value = -846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 122.708336
$ws.Range("I7").Value = 105.90909
$ws.Range("J7").Value = 136.92308
$ws.Range("K7").Value = 105.90909
$ws.Range("L7").Value = 136.92308
$ws.Range("M7").Value = 7.090909999999994
$ws.Range("N7").Value = -362.92308
$ws.Range("H16").Value = 1097.68
$ws.Range("I16").Value = 500.08334
$ws.Range("J16").Value = 1649.3077
$ws.Range("K16").Value = 500.08334
$ws.Range("L16").Value = 1649.3077
$ws.Range("M16").Value = -213.08334
$ws.Range("N16").Value = -2223.3077
$ws.Range("H31").Value = 3137.2715
$ws.Range("I31").Value = 2917.5264
$ws.Range("J31").Value = 4100.769
$ws.Range("K31").Value = 2917.5264
$ws.Range("L31").Value = 4100.769
$ws.Range("M31").Value = -2622.5264
$ws.Range("N31").Value = -4690.769
$ws.Range("H34").Value = 3137.2715
$ws.Range("I34").Value = 2917.5264
$ws.Range("J34").Value = 4100.769
$ws.Range("K34").Value = 2917.5264
$ws.Range("L34").Value = 4100.769
$ws.Range("M34").Value = -2715.5264
$ws.Range("N34").Value = -4504.769
$ws.Range("H113").Value = 1097.68
$ws.Range("I113").Value = 500.08334
$ws.Range("J113").Value = 1649.3077
$ws.Range("K113").Value = 500.08334
$ws.Range("L113").Value = 1649.3077
$ws.Range("M113").Value = 1669.91666
$ws.Range("N113").Value = -5989.3077
$ws.Range("H132").Value = 5523.7026
$ws.Range("I132").Value = 7726.278
$ws.Range("J132").Value = 3437.0527
$ws.Range("K132").Value = 23178.834
$ws.Range("L132").Value = 10311.1581
$ws.Range("M132").Value = -20648.834
$ws.Range("N132").Value = -15371.1581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 26245.264
$ws.Range("I131").Value = 1205.625
$ws.Range("J131").Value = 44455.91
$ws.Range("K131").Value = 3616.875
$ws.Range("L131").Value = 133367.73
$ws.Range("M131").Value = 1423.125
$ws.Range("N131").Value = -143447.73

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6024.558
$ws.Range("I70").Value = 5514.222
$ws.Range("J70").Value = 6392
$ws.Range("K70").Value = 5514.222
$ws.Range("L70").Value = 6392
$ws.Range("M70").Value = -5244.222
$ws.Range("N70").Value = -6932
$ws.Range("H73").Value = 6024.558
$ws.Range("I73").Value = 5514.222
$ws.Range("J73").Value = 6392
$ws.Range("K73").Value = 5514.222
$ws.Range("L73").Value = 6392
$ws.Range("M73").Value = -4578.222
$ws.Range("N73").Value = -8264
$ws.Range("H80").Value = 7755
$ws.Range("I80").Value = 12261
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 12261
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -11263
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 7755
$ws.Range("I83").Value = 12261
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 61305
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -56313
$ws.Range("N83").Value = -29984
$ws.Range("H113").Value = 1913.2354
$ws.Range("I113").Value = 1764.5
$ws.Range("J113").Value = 2045.4445
$ws.Range("K113").Value = 1764.5
$ws.Range("L113").Value = 2045.4445
$ws.Range("M113").Value = 405.5
$ws.Range("N113").Value = -6385.4445
$ws.Range("H122").Value = 2824.3635
$ws.Range("I122").Value = 2729.6296
$ws.Range("J122").Value = 3250.6667
$ws.Range("K122").Value = 8188.888800000001
$ws.Range("L122").Value = 9752.000100000001
$ws.Range("M122").Value = -5738.888800000001
$ws.Range("N122").Value = -14652.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 937.2
$ws.Range("I22").Value = 828.3333
$ws.Range("J22").Value = 1100.5
$ws.Range("K22").Value = 828.3333
$ws.Range("L22").Value = 1100.5
$ws.Range("M22").Value = -533.3333
$ws.Range("N22").Value = -1690.5
$ws.Range("H27").Value = 937.2
$ws.Range("I27").Value = 828.3333
$ws.Range("J27").Value = 1100.5
$ws.Range("K27").Value = 828.3333
$ws.Range("L27").Value = 1100.5
$ws.Range("M27").Value = -721.3333
$ws.Range("N27").Value = -1314.5
$ws.Range("H61").Value = 921460.75
$ws.Range("I61").Value = 22917.6
$ws.Range("J61").Value = 1670246.6
$ws.Range("K61").Value = 22917.6
$ws.Range("L61").Value = 1670246.6
$ws.Range("M61").Value = -22715.6
$ws.Range("N61").Value = -1670650.6
$ws.Range("H113").Value = 921460.75
$ws.Range("I113").Value = 22917.6
$ws.Range("J113").Value = 1670246.6
$ws.Range("K113").Value = 22917.6
$ws.Range("L113").Value = 1670246.6
$ws.Range("M113").Value = -20747.6
$ws.Range("N113").Value = -1674586.6
$ws.Range("H132").Value = 3784.35
$ws.Range("I132").Value = 3710.5217
$ws.Range("K132").Value = 11131.5651
$ws.Range("M132").Value = -8601.5651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1062.125
$ws.Range("I107").Value = 639.4
$ws.Range("J107").Value = 1766.6666
$ws.Range("K107").Value = 1918.2
$ws.Range("L107").Value = 5299.9998
$ws.Range("M107").Value = 1.800000000000182
$ws.Range("N107").Value = -9139.9998
$ws.Range("H132").Value = 2071.2856
$ws.Range("I132").Value = 1020.6316
$ws.Range("J132").Value = 3318.9375
$ws.Range("K132").Value = 3061.8948
$ws.Range("L132").Value = 9956.8125
$ws.Range("M132").Value = -531.8948
$ws.Range("N132").Value = -15016.8125
